$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear cells that become empty ---
$ws.Range("A13").Clear()
$ws.Range("A14").Clear()
$ws.Range("B19").Clear()
$ws.Range("C19").Clear()

# --- Set cell values to their final content ---
$ws.Range("B10").Value = "Dar conhecimentos aos alunos de noções básicas sobre ecologia e impacto das atividades da engenharia sobre o meio ambiente. Conceitos legais e institucionais para o desenvolvimento sustentável."
$ws.Range("C10").Value = "Dar conhecimentos aos alunos de noções básicas sobre ecologia e impacto das atividades da engenharia sobre o meio ambiente. Conceitos legais e institucionais para o desenvolvimento sustentável."
$ws.Range("B13").Value = "9146830 - Danúbia Caporusso Bargos"
$ws.Range("C13").Value = "9146830 - Danúbia Caporusso Bargos"
$ws.Range("B14").Value = "5464150 - Mariana Consiglio Kasemodel"
$ws.Range("C14").Value = "5464150 - Mariana Consiglio Kasemodel"
$ws.Range("A15").Value = "Programa resumido:"
$ws.Range("B15").Value = "Conceitos e Definições. Questões Ambientais. Desenvolvimento Sustentável. Desempenho Ambiental. Processos Ambientais. Norma Ambiental."
$ws.Range("C15").Value = "Conceitos e Definições. Questões Ambientais. Desenvolvimento Sustentável. Desempenho Ambiental. Processos Ambientais. Norma Ambiental."
$ws.Range("A16").Value = "Short syllabus:"
$ws.Range("B16").Value = "Concepts and Definitions. Environmental issues. Sustainable Development. Environmental performance. Environmental processes. Environmental standard."
$ws.Range("C16").Value = "Concepts and Definitions. Environmental issues. Sustainable Development. Environmental performance. Environmental processes. Environmental standard."
$ws.Range("A17").Value = "Programa:"
$ws.Range("B17").Value = "CONCEITOS E DEFINIÇÕES. Engenharia Ambiental. Meio Ambiente. Poluição Ambiental. Componentes Ambientais Críticos. QUESTÕES AMBIENTAIS. O Sujeito das Transformações Ambientais. Energia e o Meio Ambiente. Impactos Ambientais nos Três Meios. Equilíbrio Ameaçado. DESENVOLVIMENTO SUSTENTÁVEL. Conceitos Básicos. Aspectos legais. DESEMPENHO AMBIENTAL. Monitoramento Ambiental. Abrangência do Desenvolvimento Sustentável. Definição de Indicadores. Definição de Indicadores Sustentáveis. Indicadores de Desenvolvimento Humano  IDH. Indicadores de Sustentabilidade Ambiental. Controle de Processos Ambientais. PROCESSOS AMBIENTAIS. Controle Processo ETA. Água na Natureza. Caracterização da Água. Indicadores de Qualidade da Água. NORMA AMBIENTAL. Portaria 518. CONAMA 20. Desastre Ecológico."
$ws.Range("C17").Value = "CONCEITOS E DEFINIÇÕES. Engenharia Ambiental. Meio Ambiente. Poluição Ambiental. Componentes Ambientais Críticos. QUESTÕES AMBIENTAIS. O Sujeito das Transformações Ambientais. Energia e o Meio Ambiente. Impactos Ambientais nos Três Meios. Equilíbrio Ameaçado. DESENVOLVIMENTO SUSTENTÁVEL. Conceitos Básicos. Aspectos legais. DESEMPENHO AMBIENTAL. Monitoramento Ambiental. Abrangência do Desenvolvimento Sustentável. Definição de Indicadores. Definição de Indicadores Sustentáveis. Indicadores de Desenvolvimento Humano  IDH. Indicadores de Sustentabilidade Ambiental. Controle de Processos Ambientais. PROCESSOS AMBIENTAIS. Controle Processo ETA. Água na Natureza. Caracterização da Água. Indicadores de Qualidade da Água. NORMA AMBIENTAL. Portaria 518. CONAMA 20. Desastre Ecológico."
$ws.Range("A18").Value = "Syllabus:"
$ws.Range("B18").Value = "CONCEPTS AND DEFINITIONS. Environmental Engineering. Environment. Environmental pollution. Environmental Critical Components. ENVIRONMENTAL ISSUES. The Subject of Environmental Transformations. Energy and the Environment. Environmental impacts in the three media. Threatened equilibrium. SUSTAINABLE DEVELOPMENT. Basic Concepts. Legal aspects. ENVIRONMENTAL PERFORMANCE. Environmental Monitoring. Scope of Sustainable Development. Definition of indicators. Definition of Sustainable Indicators. Human Development Indicators - HDI. Environmental Sustainability Indicators. Control of Environmental Processes. ENVIRONMENTAL PROCESSES. ETA Process Control. Water in Nature. Characterization of Water. Water Quality Indicators. ENVIRONMENTAL STANDARD. Ordinance 518. CONAMA 20. Ecological disaster."
$ws.Range("C18").Value = "CONCEPTS AND DEFINITIONS. Environmental Engineering. Environment. Environmental pollution. Environmental Critical Components. ENVIRONMENTAL ISSUES. The Subject of Environmental Transformations. Energy and the Environment. Environmental impacts in the three media. Threatened equilibrium. SUSTAINABLE DEVELOPMENT. Basic Concepts. Legal aspects. ENVIRONMENTAL PERFORMANCE. Environmental Monitoring. Scope of Sustainable Development. Definition of indicators. Definition of Sustainable Indicators. Human Development Indicators - HDI. Environmental Sustainability Indicators. Control of Environmental Processes. ENVIRONMENTAL PROCESSES. ETA Process Control. Water in Nature. Characterization of Water. Water Quality Indicators. ENVIRONMENTAL STANDARD. Ordinance 518. CONAMA 20. Ecological disaster."
$ws.Range("A19").Value = "Avaliação:"
$ws.Range("A20").Value = "Método:"
$ws.Range("B20").Value = "Aulas expositivas com a utilização de recursos de projeções e audiovisual."
$ws.Range("C20").Value = "Aulas expositivas com a utilização de recursos de projeções e audiovisual."
$ws.Range("A21").Value = "Critério:"
$ws.Range("B21").Value = "Média ponderada de 2 avaliações escritas com nota final (NF ≥ 5,0)"
$ws.Range("C21").Value = "Média ponderada de 2 avaliações escritas com nota final (NF ≥ 5,0)"
$ws.Range("A22").Value = "Norma de recuperação:"
$ws.Range("B22").Value = "(NF+RC)/2 ≥ 5,0, onde RC é uma prova escrita de recuperação a ser aplicada"
$ws.Range("C22").Value = "(NF+RC)/2 ≥ 5,0, onde RC é uma prova escrita de recuperação a ser aplicada"
$ws.Range("A23").Value = "Bibliografia:"
$ws.Range("B23").Value = "1)        BRAGA, B.; HESPANHOL, I.; CONEJO, J. G. L.; MIERZWA, J. C.; BARROS, M. T. L.; SPENCER, M.; PORTO, M.; NUCCI, N.; JULIANO, N.; EIGER, S. Introdução à Engenharia Ambiental: O Desafio do Desenvolvimento Sustentável. Pearson (2ª Edição), 336 p., 2005.2)        VESILIND, P.A.; MORGAN, S. M.; HEINE, L. G. Introdução à Engenharia Ambiental. Cengage (3ª edição), 472 p., 2018.3)        CALIJURI, M. C.; CUNHA, D. G. F. Engenharia Ambiental: Conceitos, Tecnologias e Gestão. Elsevier (1ª Edição), 832 p., 2012.4)        CAPAZ, R. S.; HORTA NOGUEIRA, L. A. Ciências Ambientais para Engenharia. Elsevier (1ª Edição), 252 p., 2014.5)        DAVIS, M. L.; MASTEN, S. J. Princípios de Engenharia Ambiental. Mc Graw Hill Educations (3ª Edição), 872 p., 2016;"
$ws.Range("C23").Value = "1)        BRAGA, B.; HESPANHOL, I.; CONEJO, J. G. L.; MIERZWA, J. C.; BARROS, M. T. L.; SPENCER, M.; PORTO, M.; NUCCI, N.; JULIANO, N.; EIGER, S. Introdução à Engenharia Ambiental: O Desafio do Desenvolvimento Sustentável. Pearson (2ª Edição), 336 p., 2005.2)        VESILIND, P.A.; MORGAN, S. M.; HEINE, L. G. Introdução à Engenharia Ambiental. Cengage (3ª edição), 472 p., 2018.3)        CALIJURI, M. C.; CUNHA, D. G. F. Engenharia Ambiental: Conceitos, Tecnologias e Gestão. Elsevier (1ª Edição), 832 p., 2012.4)        CAPAZ, R. S.; HORTA NOGUEIRA, L. A. Ciências Ambientais para Engenharia. Elsevier (1ª Edição), 252 p., 2014.5)        DAVIS, M. L.; MASTEN, S. J. Princípios de Engenharia Ambiental. Mc Graw Hill Educations (3ª Edição), 872 p., 2016;"

# --- Adjust row heights ---
$ws.Rows.Item(13).AutoFit()
$ws.Rows.Item(14).AutoFit()
$ws.Rows.Item(15).RowHeight = 60
$ws.Rows.Item(16).RowHeight = 60
$ws.Rows.Item(17).RowHeight = 120
$ws.Rows.Item(18).RowHeight = 120
$ws.Rows.Item(19).AutoFit()
$ws.Rows.Item(21).RowHeight = 60
$ws.Rows.Item(22).RowHeight = 60
$ws.Rows.Item(23).RowHeight = 120